$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ============================================================
# 1) Add new column O (year 2023) by cloning column N formatting
# ============================================================
for ($r = 2; $r -le 33; $r++) {
    $ws.Range("N" + $r).Copy()
    $ws.Range("O" + $r).PasteSpecial(-4122)
}

# Fix up the handful of O-column rows whose target style differs from N
$ws.Range("N8").Copy()
$ws.Range("O15").PasteSpecial(-4122)
$ws.Range("N8").Copy()
$ws.Range("O17").PasteSpecial(-4122)
$ws.Range("N8").Copy()
$ws.Range("O23").PasteSpecial(-4122)

# O13 needs a brand-new style: same as N13 but with right-aligned horizontal
$ws.Range("N13").Copy()
$ws.Range("O13").PasteSpecial(-4122)
$ws.Range("O13").HorizontalAlignment = -4152

# Write the actual 2023 values into column O
$ws.Range("O3").Value = 2023
$ws.Range("O4").Value = 2.3944505088207331
$ws.Range("O5").Value = 0.52951886522435021
$ws.Range("O6").Value = 4.3000500913119915
$ws.Range("O7").Value = 0.51981356020307379
$ws.Range("O8").Value = "-"
$ws.Range("O9").Value = 1.028739554007112
$ws.Range("O10").Value = 0.45336944169064486
$ws.Range("O11").Value = 0.15211161341746121
$ws.Range("O12").Value = 0.75073722395392273
$ws.Range("O13").Value = "-"
$ws.Range("O14").Value = "-"
$ws.Range("O15").Value = "-"
$ws.Range("O16").Value = 0.32236434908190637
$ws.Range("O17").Value = "-"
$ws.Range("O18").Value = 0.63756806039044667
$ws.Range("O19").Value = 0.47449906455898705
$ws.Range("O20").Value = 0.13658701822343999
$ws.Range("O21").Value = 0.80742182138214469
$ws.Range("O22").Value = 1.4528231986808364
$ws.Range("O23").Value = "-"
$ws.Range("O24").Value = 2.8788783889796536
$ws.Range("O25").Value = 7.2371521991664292
$ws.Range("O26").Value = 0.1843916182945988
$ws.Range("O27").Value = 14.38048139128356
$ws.Range("O28").Value = 5.7129502194292243
$ws.Range("O29").Value = 2.3720702955125175
$ws.Range("O30").Value = 9.7530937387050578
$ws.Range("O31").Value = 1.3736037318066185
$ws.Range("O32").Value = 0.56245500359971201
$ws.Range("O33").Value = 2.148066203400389

# ============================================================
# 2) Row height / layout tweaks on row 1
# ============================================================
$ws.Rows.Item(1).RowHeight = 42.75

# ============================================================
# 3) Capitalise the gender labels (aialdar/erkekter etc.)
#    "women" rows: A/B/C -> Аялдар / Женщины / Woman
#    "men" rows:   A/B/C -> Эркектер / Мужчины / Men
# ============================================================
$ws.Range("A5").Value = "Аялдар"
$ws.Range("B5").Value = "Женщины"
$ws.Range("C5").Value = "Woman"
$ws.Range("A8").Value = "Аялдар"
$ws.Range("B8").Value = "Женщины"
$ws.Range("C8").Value = "Woman"
$ws.Range("A11").Value = "Аялдар"
$ws.Range("B11").Value = "Женщины"
$ws.Range("C11").Value = "Woman"
$ws.Range("A14").Value = "Аялдар"
$ws.Range("B14").Value = "Женщины"
$ws.Range("C14").Value = "Woman"
$ws.Range("A17").Value = "Аялдар"
$ws.Range("B17").Value = "Женщины"
$ws.Range("C17").Value = "Woman"
$ws.Range("A20").Value = "Аялдар"
$ws.Range("B20").Value = "Женщины"
$ws.Range("C20").Value = "Woman"
$ws.Range("A23").Value = "Аялдар"
$ws.Range("B23").Value = "Женщины"
$ws.Range("C23").Value = "Woman"
$ws.Range("A26").Value = "Аялдар"
$ws.Range("B26").Value = "Женщины"
$ws.Range("C26").Value = "Woman"
$ws.Range("A29").Value = "Аялдар"
$ws.Range("B29").Value = "Женщины"
$ws.Range("C29").Value = "Woman"
$ws.Range("A32").Value = "Аялдар"
$ws.Range("B32").Value = "Женщины"
$ws.Range("C32").Value = "Woman"
$ws.Range("A6").Value = "Эркектер"
$ws.Range("B6").Value = "Мужчины"
$ws.Range("C6").Value = "Men"
$ws.Range("A9").Value = "Эркектер"
$ws.Range("B9").Value = "Мужчины"
$ws.Range("C9").Value = "Men"
$ws.Range("A12").Value = "Эркектер"
$ws.Range("B12").Value = "Мужчины"
$ws.Range("C12").Value = "Men"
$ws.Range("A15").Value = "Эркектер"
$ws.Range("B15").Value = "Мужчины"
$ws.Range("C15").Value = "Men"
$ws.Range("A18").Value = "Эркектер"
$ws.Range("B18").Value = "Мужчины"
$ws.Range("C18").Value = "Men"
$ws.Range("A21").Value = "Эркектер"
$ws.Range("B21").Value = "Мужчины"
$ws.Range("C21").Value = "Men"
$ws.Range("A24").Value = "Эркектер"
$ws.Range("B24").Value = "Мужчины"
$ws.Range("C24").Value = "Men"
$ws.Range("A27").Value = "Эркектер"
$ws.Range("B27").Value = "Мужчины"
$ws.Range("C27").Value = "Men"
$ws.Range("A30").Value = "Эркектер"
$ws.Range("B30").Value = "Мужчины"
$ws.Range("C30").Value = "Men"
$ws.Range("A33").Value = "Эркектер"
$ws.Range("B33").Value = "Мужчины"
$ws.Range("C33").Value = "Men"

# ============================================================
# 4) Reset the stray stale selection left over from editing
# ============================================================
$ws.Range("A1").Select()
